# UNFCCC_scaling_mapping.xlsx update
#
# - "method" sheet: collapse the per-sector (energy / industry /
#   domestic_transport / buildings) constant/linear extension rows down to
#   a single generic "NA" placeholder row.
# - "year" sheet: same collapse, down to a single "NA" placeholder row.
# - "map" sheet is unaffected content-wise; it just becomes the active tab.
#
# Removing the now-unused "constant" / "linear" / "all" shared strings and
# adding "NA" happens automatically when the workbook is written back out,
# since those strings are no longer referenced after the edits below.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # map
$ws2 = $wb.Worksheets.Item(2)   # method
$ws3 = $wb.Worksheets.Item(3)   # year

# --- "method" sheet --------------------------------------------------
# Replace every data row with a single "NA" row, then drop the rest.
$ws2.Range("A2:E2").Value = "NA"
$ws2.Rows("3:5").Delete()
$ws2.Range("B2").ClearFormats()
$ws2.Range("D2:E2").Select() | Out-Null

# --- "year" sheet -----------------------------------------------------
$ws3.Range("A2:D2").Value = "NA"
$ws3.Rows("3:5").Delete()
$ws3.Range("B2").ClearFormats()
$ws3.Range("A2:D2").Select() | Out-Null

# --- make "map" the active / selected sheet ---------------------------
$ws1.Activate()
